$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns
# to the latest scraped values. Price strings that look like plain
# numbers are forced to remain text (matching the original inline
# string storage) and then have their format reset so no stray
# cell style is introduced.

$ws.Range('D2').Value = '67.100.24'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').Value = '3.114.25'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.08'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '173.72'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.35%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.110.52'
$ws.Range('E8').Value = '  +2.80%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.45'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.34%  '
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.483'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.23'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D16').Value = '3.629.40'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '67.042.92'
$ws.Range('E17').Value = '  +1.31%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '7.19'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').Value = '3.115.24'
$ws.Range('E19').Value = '  +2.81%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.20'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '486.68'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.86%  '
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.60'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '84.40'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.37'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +4.36%  '
$ws.Range('E26').Value = '  +4.19%  '
$ws.Range('E27').Value = '  +0.52%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.03'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.40'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.88%  '
$ws.Range('E31').Value = '  +2.52%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '28.87'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.14%  '
$ws.Range('E33').Value = '  +0.86%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.114'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.10%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.989'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.45%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '47.72'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('E39').Value = '  +3.22%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '50.19'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.11%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.315'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.96%  '
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('E43').Value = '  +0.80%  '
$ws.Range('E44').Value = '  -0.96%  '
$ws.Range('D45').Value = '2.845.62'
$ws.Range('E45').Value = '  +4.60%  '
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '383.93'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '137.08'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.04%  '
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '25.16'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('E51').Value = '  +0.46%  '
